# Applies the "Time Dilation" -> "Chemistry" rewrite to the active document.
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $found) {
        Write-Host "NOT FOUND:" $old
    }
}

# --- Title & byline -------------------------------------------------------
Replace-Text "Time Dilation: A Twist in Spacetime" "The Beauty of Separation: Uncovering the Wonders of Chemistry"
Replace-Text "Isaac Newton" "Lucy Phillips"

# --- Email line: collapse "isaac" + "." + "newton@physics" runs into one --
Replace-Text "isaac.newton@physics" "lphillips123@educonnect"

# --- Body paragraph 1 (physics intro -> chemistry intro) ------------------
Replace-Text "In the realm of physics, where time and space intertwine, lies a fascinating phenomenon that challenges our perception of reality: time dilation" "Immerse yourselves in the captivating realm of chemistry, where matter undergoes transformations, revealing the intricacies of the world around us"

Replace-Text " This remarkable effect arises from the interplay between the speed of light and the curvature of spacetime, unveiling a universe where time flows differently for different observers" " Chemistry, a branch of science that delves into the fundamental principles governing the interactions between substances, offers a gateway to understanding the composition and behavior of matter, unveiling the building blocks of life itself"

Replace-Text " As we journey through the cosmos, the faster we travel, the slower time passes for us, a concept that has profound implications for our understanding of the universe and our place within it" " From the intricate dance of atoms in molecular structures to the vibrant interplay of elements in chemical reactions, chemistry holds the key to unlocking the secrets of the material world"

# --- Gravity paragraph -----------------------------------------------------
Replace-Text "Gravity, the invisible force that binds us to the Earth and governs the motion of celestial bodies, also plays a crucial role in shaping the fabric of spacetime" "Chemistry unveils the microscopic world, revealing the subatomic particles that orchestrate the chemical symphony: electrons, protons, and neutrons"

Replace-Text " The presence of massive objects, such as planets, stars, and black holes, warps spacetime, creating regions where time elapses at different rates" " Their intricate interactions, driven by fundamental forces, govern the properties and behaviors of elements, defining their place in the periodic table, the compass of chemistry"

Replace-Text " As we venture closer to these gravitational behemoths, time slows down, leading to remarkable effects that have been experimentally verified and continue to captivate the minds of scientists and philosophers alike" " As we delve deeper into the realms of chemical reactions, the stage is set for an awe-inspiring spectacle: substances transforming into new substances, accompanied by energy exchanges, shaping the very fabric of our world"

# --- "Time dilation ... revolution" paragraph -----------------------------
Replace-Text "Time dilation, a consequence of Einstein's Theory of Special Relativity, has spurred a revolution in our comprehension of the universe" "The significance of chemistry extends far beyond the confines of the laboratory; it infiltrates every aspect of our lives"

Replace-Text " From the intricacies of black hole physics to the mind-boggling implications of interstellar travel, time dilation challenges our conventional notions of time and space and opens up a realm of possibilities that were once thought to be beyond our reach" " From the nourishment we derive from food, the healing power of medicines, to the materials that shape our technological world, chemistry plays an indispensable role"

# Two brand-new sentences are appended after that last run (before the
# paragraph-closing "." run), so splice them in right before it. The search
# text is long/unique enough to only match this one spot in the document.
$rng = $d.Content
$find2 = $rng.Find
$find2.ClearFormatting()
$find2.Text = "the materials that shape our technological world, chemistry plays an indispensable role"
if ($find2.Execute()) {
    $rng.Collapse(0)
    $rng.InsertAfter(". It is a field that constantly evolves, revealing new insights and applications, driving innovation and propelling humanity towards a future filled with endless possibilities")
} else {
    Write-Host "SPLICE POINT 1 NOT FOUND"
}

# --- Summary heading paragraph ---------------------------------------------
Replace-Text "Time dilation, a product of the interplay between the speed of light and the curvature of spacetime, is a profound phenomenon that alters our perception of time" "Chemistry, the study of matter and its interactions, uncovers the fundamental principles governing the composition and behavior of matter"

Replace-Text " As we approach the speed of light or venture closer to massive objects, time slows down, a concept that has been experimentally verified and corroborated by various observations" " Delving into the microscopic realm of atoms and molecules, chemistry reveals the intricate dance of subatomic particles that shape the properties of elements"

Replace-Text " This remarkable effect has profound implications for our understanding of the universe, gravity, and the nature of time itself, inspiring awe and wonder in our quest to unravel the mysteries of the cosmos" " Unveiling the secrets of chemical reactions, it showcases the transformative power of substances as they morph into new entities, accompanied by energy exchanges"

$rng2 = $d.Content
$find4 = $rng2.Find
$find4.ClearFormatting()
$find4.Text = "morph into new entities, accompanied by energy exchanges"
if ($find4.Execute()) {
    $rng2.Collapse(0)
    $rng2.InsertAfter(". Chemistry's significance extends beyond the laboratory, permeating our daily lives through food, medicines, and technological marvels. As a dynamic field constantly evolving, chemistry holds the promise of new discoveries and applications, driving innovation and shaping the future of humanity")
} else {
    Write-Host "SPLICE POINT 2 NOT FOUND"
}

# --- Trailing empty paragraph added at the very end of the body -----------
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.Text = "`r"
